$wb = $excel.ActiveWorkbook

# Rename worksheets
$wb.Worksheets.Item(1).Name = "GNG_TO-16509961167546022"
$wb.Worksheets.Item(2).Name = "NB_TO-1650996120474577"
$wb.Worksheets.Item(3).Name = "RS_TO-1650996120474577"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509961205385995"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509961206026"

# Sheet 1 (GNG) updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509961167225628.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961167385623.csv"
$ws1.Range("B4").Value = "go_stims-16509961167385623.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961167546022.csv"

# Sheet 2 (NB) updates
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_2-16509961178265638.csv"
$ws2.Range("B3").Value = "TB-16509961201785607.csv"
$ws2.Range("B4").Value = "ZB-match_2-16509961177706017.csv"
$ws2.Range("B5").Value = "ZB-match_5-16509961175945609.csv"
$ws2.Range("B6").Value = "OB-1650996118050599.csv"
$ws2.Range("B7").Value = "OB-16509961186106014.csv"
$ws2.Range("B8").Value = "TB-165099612045056.csv"
$ws2.Range("B9").Value = "OB-1650996118322566.csv"
$ws2.Range("B10").Value = "TB-16509961198586001.csv"

# Sheet 4 (TOL) updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509961205065947.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961204825668.csv"
$ws4.Range("B4").Value = "MM_stims-16509961205225663.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961205065947.csv"
$ws4.Range("B6").Value = "MM_stims-16509961205385995.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961205225663.csv"

# Sheet 5 (vSAT) updates
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16509961205865636.csv"
$ws5.Range("B3").Value = "SAT_stims-16509961205546007.csv"
$ws5.Range("B4").Value = "SAT_stims-16509961205385995.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650996120570603.csv"
